# PrecioFrutaHortalizas - Vega Monumental Concepción - Espinaca
# Weekly update: insert a new weekly record as row 124, pushing the
# existing historical rows (124-149) down by one (125-150).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 124, shifting rows
# 124:149 down to 125:150 (mirrors Excel's Rows().Insert()).
$ws.Rows.Item(124).Insert()

# Populate the newly inserted row 124 with the new weekly price record.
$ws.Range("A124").Value = 11
$ws.Range("B124").Value = "Vega Monumental Concepción"
$ws.Range("C124").Value = "Bíobío"
$ws.Range("D124").Value = 45258
$ws.Range("E124").Value = 8
$ws.Range("F124").Value = 100112012
$ws.Range("G124").Value = "Espinaca"
$ws.Range("H124").Value = "Sin especificar"
$ws.Range("I124").Value = "Primera"
$ws.Range("J124").Value = 150
$ws.Range("K124").Value = 14000
$ws.Range("L124").Value = 14000
$ws.Range("M124").Value = 14000
$ws.Range("N124").Value = "`$/cuna 10 kilos"
$ws.Range("O124").Value = "Región Metropolitana"
$ws.Range("P124").Value = 1400
$ws.Range("Q124").Value = 10
$ws.Range("R124").Value = "Hortaliza"
